# Common: Added some boirelplate stuff
# Append new "lab.mixture.*" / "lab.vape.*" translation rows to the
# "Import" sheet of the translations workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# language code used in column A for every data row
$lang = "cs"

# Row data: row number, label (col B), translation (col C)
$rows = @(
    @(301, "lab.mixture.menu", "Mix"),
    @(302, "lab.mixture.title", "Mixy"),
    @(303, "lab.mixture.subtitle", "Každý požitek potřebuje liquid a tato sekce slouží pro správu namíchaných liquidů (včetně hotových); tyto mixy se pak dále používají ve vapování pro trasování, jak který mix chutnal."),
    @(304, "lab.mixture.button.create", "Nový mix"),
    @(305, "lab.mixture.button.list", "Seznam mixů"),
    @(306, "lab.mixture.create.title", "Nový mix"),
    @(307, "lab.mixture.create.subtitle", "Mix vám pomůže sledovat složení liquidu, množství nikotinu a jeho stáří (tzn. zrání)."),
    @(308, "lab.mixture.list.title", "Seznam mixů"),
    @(309, "lab.vape.title", "Vape"),
    @(310, "lab.vape.subtitle", "Toto je nejzajímavější část aplikace, kde si můžete trasovat zážitky z vapování, které vám tak umožní zjistit, jaké spirálky, nastavení vzduchu, liquidy (a jejich stáří) vám vyhovují nejvíce."),
    @(311, "lab.vape.button.create", "Nový vape"),
    @(312, "lab.vape.button.list", "Seznam vapů"),
    @(313, "lab.vape.create.title", "Nový vape"),
    @(314, "lab.vape.create.subtitle", "Vape je základní stavební kámen pro záznam chutě a požitku z vapování."),
    @(315, "lab.vape.list.title", "Seznam vapů")
)

# Rows whose translation text is long enough that Excel's wrap-text
# column formatting auto-expands the row height (matches existing rows
# like 264/266 which already carry ht="39" for similarly long strings).
$tallRows = @(303, 310)

foreach ($row in $rows) {
    $r = $row[0]
    $label = $row[1]
    $translation = $row[2]

    $ws.Cells.Item($r, 1).Value2 = $lang
    $ws.Cells.Item($r, 2).Value2 = $label
    $ws.Cells.Item($r, 3).Value2 = $translation

    # Copy the formatting (style index, incl. wrapText) from the last
    # existing data row so the new rows match the established look.
    $ws.Range("A300:C300").Copy()
    $ws.Range("A$r" + ":C$r").PasteSpecial(-4122)

    if ($tallRows -contains $r) {
        $ws.Rows.Item($r).RowHeight = 39
    }
}

$excel.CutCopyMode = 0

# Update the current selection to mirror what a human editor typing
# these rows in the UI would leave behind.
$null = $ws.Range("B310").Select()

"Added $($rows.Count) translation rows (301-315) to Import sheet"
